# Updated symbol list (price / volume(1h) refresh) for cryptos.xlsx
# Cells hold their numbers/percentages as literal text (t="inlineStr" in the
# source workbook), so values are written with a leading apostrophe to force
# Excel to keep them as text instead of auto-converting to numeric/percent
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.97"
$ws.Range("E2").Value = "'0.98%"
$ws.Range("D3").Value = "'32.08"
$ws.Range("E3").Value = "'1.12%"
$ws.Range("D4").Value = "'4.989"
$ws.Range("E4").Value = "'-2.27%"
$ws.Range("D5").Value = "'0.07913"
$ws.Range("E5").Value = "'-3.62%"
$ws.Range("D6").Value = "'2.104"
$ws.Range("E6").Value = "'-18.35%"
$ws.Range("D7").Value = "'7.855"
$ws.Range("E7").Value = "'0.15%"
$ws.Range("D8").Value = "'3.797"
$ws.Range("E8").Value = "'-1.15%"
$ws.Range("D9").Value = "'0.9286"
$ws.Range("E9").Value = "'-0.10%"
$ws.Range("D10").Value = "'0.1752"
$ws.Range("E10").Value = "'-0.42%"
$ws.Range("D11").Value = "'0.08041"
$ws.Range("E11").Value = "'7.41%"
$ws.Range("D12").Value = "'0.08818"
$ws.Range("E12").Value = "'-1.11%"
$ws.Range("D13").Value = "'0.03131"
$ws.Range("E13").Value = "'4.16%"
$ws.Range("E14").Value = "'0.25%"
$ws.Range("D15").Value = "'0.001543"
$ws.Range("E15").Value = "'1.55%"
$ws.Range("D16").Value = "'0.005927"
$ws.Range("E16").Value = "'0.38%"
$ws.Range("E17").Value = "'-3.55%"
$ws.Range("E19").Value = "'1.51%"
$ws.Range("D20").Value = "'0.1290"
$ws.Range("E20").Value = "'-3.46%"
$ws.Range("D21").Value = "'4.157"
$ws.Range("E21").Value = "'6.35%"
$ws.Range("E23").Value = "'-0.17%"
$ws.Range("D24").Value = "'0.001236"
$ws.Range("E24").Value = "'-0.79%"
$ws.Range("D25").Value = "'0.004508"
$ws.Range("E25").Value = "'-1.09%"
$ws.Range("D26").Value = "'0.0001248"
$ws.Range("E26").Value = "'4.34%"
$ws.Range("D39").Value = "'0.01733"
$ws.Range("E39").Value = "'-2.34%"
$ws.Range("D40").Value = "'0.04949"
$ws.Range("E40").Value = "'8.15%"
$ws.Range("D41").Value = "'0.007379"
$ws.Range("E41").Value = "'7.52%"
$ws.Range("D42").Value = "'0.1368"
$ws.Range("E42").Value = "'-1.00%"
$ws.Range("D43").Value = "'0.002307"
$ws.Range("E43").Value = "'4.70%"
$ws.Range("E44").Value = "'15.51%"
$ws.Range("D45").Value = "'0.00006064"
$ws.Range("E45").Value = "'-2.39%"
$ws.Range("E46").Value = "'0.23%"
$ws.Range("D48").Value = "'0.8234"
$ws.Range("E48").Value = "'1.31%"
$ws.Range("E49").Value = "'0.23%"
$ws.Range("E50").Value = "'0.23%"
